# Generate Report for Handback
# Replaces the two source file identifiers (and their dependent handoff/handback
# xliff names + timestamps) across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "d500ff09-c503-46ef-896e-990d401aa645.md"
$wsOverview.Range("B2").Value = "e2e\d500ff09-c503-46ef-896e-990d401aa645.md"
$wsOverview.Range("G2").Value = "2016-09-06 19:24:42"

$wsOverview.Range("A3").Value = "ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
$wsOverview.Range("B3").Value = "e2e\ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
$wsOverview.Range("G3").Value = "2016-09-06 19:24:42"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq "e2e\820a4fd4-c8a5-43eb-b1d4-d2ffb3d3408b.md") {
        $hl.TextToDisplay = "e2e\d500ff09-c503-46ef-896e-990d401aa645.md"
    }
    elseif ($hl.TextToDisplay -eq "e2e\88e1953d-2ee2-4956-8a33-b5a046a4c563.md") {
        $hl.TextToDisplay = "e2e\ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "d500ff09-c503-46ef-896e-990d401aa645.md"
$wsZhCn.Range("G2").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-06 19:24:37"
$wsZhCn.Range("I2").Value = "d500ff09-c503-46ef-896e-990d401aa645.md"
$wsZhCn.Range("J2").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-06 19:24:56"

$wsZhCn.Range("A3").Value = "ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
$wsZhCn.Range("G3").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-06 19:24:37"
$wsZhCn.Range("I3").Value = "ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
$wsZhCn.Range("J3").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-06 19:24:56"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.TextToDisplay -eq "820a4fd4-c8a5-43eb-b1d4-d2ffb3d3408b.md") {
        $hl.TextToDisplay = "d500ff09-c503-46ef-896e-990d401aa645.md"
    }
    elseif ($hl.TextToDisplay -eq "88e1953d-2ee2-4956-8a33-b5a046a4c563.md") {
        $hl.TextToDisplay = "ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "d500ff09-c503-46ef-896e-990d401aa645.md"
$wsDeDe.Range("G2").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-06 19:24:42"
$wsDeDe.Range("I2").Value = "d500ff09-c503-46ef-896e-990d401aa645.md"
$wsDeDe.Range("J2").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-06 19:25:23"

$wsDeDe.Range("A3").Value = "ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
$wsDeDe.Range("G3").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-06 19:24:42"
$wsDeDe.Range("I3").Value = "ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
$wsDeDe.Range("J3").Value = "d500ff09-c503-46ef-896e-990d401aa645.567b6d8be8fcaaa590f5c6d77b2b4ec3ca8adac6.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-06 19:25:23"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq "820a4fd4-c8a5-43eb-b1d4-d2ffb3d3408b.md") {
        $hl.TextToDisplay = "d500ff09-c503-46ef-896e-990d401aa645.md"
    }
    elseif ($hl.TextToDisplay -eq "88e1953d-2ee2-4956-8a33-b5a046a4c563.md") {
        $hl.TextToDisplay = "ffff8de27131-a65e-4a17-a8de-6d356a985e55.md"
    }
}
